$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split the name line "AYOMIDE ESTHER, AINA" into "AYOMIDE" + ", AINA",
#    moving the hidden _GoBack bookmark to the split point (between the
#    two runs), mirroring Word's own "last edit" bookmark tracking.
# ---------------------------------------------------------------------
$nameRange = $d.Paragraphs.Item(2).Range
$nameStart = $nameRange.Start

# Offset right after "AYOMIDE" (7 characters) -> this is where the
# bookmark must sit once " ESTHER" is removed.
$splitPos = $nameStart + 7
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove " ESTHER" (7 chars right after "AYOMIDE"), leaving
# "AYOMIDE" + bookmark + ", AINA"
$deleteRange = $d.Range($splitPos, $splitPos + 7)
$deleteRange.Text = ""

# ---------------------------------------------------------------------
# 2) Remove the VOLUNTEER EXPERIENCE section in its entirety, including
#    the blank ListParagraph spacer right before it.
# ---------------------------------------------------------------------
$sectionStartPara = $null
$sectionEndPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "VOLUNTEER EXPERIENCE") {
        $sectionStartPara = $i - 1   # include the preceding blank ListParagraph
    }
    if ($txt -eq "SKILLS") {
        $sectionEndPara = $i - 1    # last paragraph of the section to remove
        break
    }
}

if ($sectionStartPara -ne $null -and $sectionEndPara -ne $null) {
    $delStart = $d.Paragraphs.Item($sectionStartPara).Range.Start
    $delEnd = $d.Paragraphs.Item($sectionEndPara).Range.End
    $d.Range($delStart, $delEnd).Delete() | Out-Null
}

# ---------------------------------------------------------------------
# 3) Drop the stale rendering hint on the HOBBIES paragraph
#    (<w:lastRenderedPageBreak/>) by re-typing its text, which forces a
#    fresh run without that cached layout artifact.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i).Range
    if ($p.Text.Trim() -eq "HOBBIES") {
        $hobbiesText = $p.Text
        $scope = $d.Range($p.Start, $p.End)
        $scope.Find.Execute($hobbiesText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $hobbiesText, 2) | Out-Null
        break
    }
}
